$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.889.22'
$ws.Range("E2").Value = '  +1.59%  '
$ws.Range("D3").Value = '2.091.76'
$ws.Range("E3").Value = '  +1.22%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("Z1").Value = '''233.25'
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("Z1").Value = '''0.625'
$ws.Range("Z1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("Z1").Value = '''57.57'
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E8").Value = '  +1.30%  '
$ws.Range("Z1").Value = '''0.389'
$ws.Range("Z1").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("Z1").Value = '''0.0782'
$ws.Range("Z1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E10").Value = '  +2.60%  '
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").Value = '2.385.84'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("Z1").Value = '''14.41'
$ws.Range("Z1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E13").Value = '  -1.33%  '
$ws.Range("Z1").Value = '''21.16'
$ws.Range("Z1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E14").Value = '  +2.49%  '
$ws.Range("Z1").Value = '''0.765'
$ws.Range("Z1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E15").Value = '  -1.41%  '
$ws.Range("Z1").Value = '''5.24'
$ws.Range("Z1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E16").Value = '  +2.65%  '
$ws.Range("D17").Value = '2.080.59'
$ws.Range("D18").Value = '37.838.48'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("Z1").Value = '''6.13'
$ws.Range("Z1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E19").Value = '  -2.35%  '
$ws.Range("Z1").Value = '''70.86'
$ws.Range("Z1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = '0.0₃0821'
$ws.Range("E21").Value = '  +1.17%  '
$ws.Range("Z1").Value = '''228.45'
$ws.Range("Z1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("Z1").Value = '''170.84'
$ws.Range("Z1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E26").Value = '  +1.99%  '
$ws.Range("E27").Value = '  +10.20%  '
$ws.Range("Z1").Value = '''8.97'
$ws.Range("Z1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E28").Value = '  +2.38%  '
$ws.Range("E29").Value = '  +0.54%  '
$ws.Range("Z1").Value = '''19.50'
$ws.Range("Z1").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E30").Value = '  +2.45%  '
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("Z1").Value = '''4.63'
$ws.Range("Z1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("Z1").Value = '''0.0627'
$ws.Range("Z1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E33").Value = '  +1.86%  '
$ws.Range("Z1").Value = '''4.61'
$ws.Range("Z1").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E34").Value = '  +0.98%  '
$ws.Range("Z1").Value = '''2.51'
$ws.Range("Z1").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("E36").Value = '  +3.91%  '
$ws.Range("Z1").Value = '''3.40'
$ws.Range("Z1").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E37").Value = '  +5.42%  '
$ws.Range("Z1").Value = '''1.00'
$ws.Range("Z1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("Z1").Value = '''5.45'
$ws.Range("Z1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E39").Value = '  -3.44%  '
$ws.Range("Z1").Value = '''0.101'
$ws.Range("Z1").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E40").Value = '  +6.94%  '
$ws.Range("E41").Value = '  -0.62%  '
$ws.Range("Z1").Value = '''97.40'
$ws.Range("Z1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("E43").Value = '  +0.80%  '
$ws.Range("D44").Value = '1.453.83'
$ws.Range("E44").Value = '  -1.12%  '
$ws.Range("Z1").Value = '''1.16'
$ws.Range("Z1").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("Z1").Value = '''1.06'
$ws.Range("Z1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E46").Value = '  +3.68%  '
$ws.Range("Z1").Value = '''15.76'
$ws.Range("Z1").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E47").Value = '  +5.07%  '
$ws.Range("Z1").Value = '''4.05'
$ws.Range("Z1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E48").Value = '  -6.32%  '
$ws.Range("Z1").Value = '''7.40'
$ws.Range("Z1").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("Z1").Value = '''3.01'
$ws.Range("Z1").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("D51").Value = '2.280.76'
$ws.Range("E51").Value = '  +0.92%  '

$excel.CutCopyMode = 0

